# "Generate Report for Handback"
# Updates the localization-status workbook to reflect a new handback run:
#  - Overview status ("Ready for handoff" -> "Handed back: in sync with en-US")
#  - zh-cn / de-de "Latest Handback DateTime" bumped to the new handback timestamps
#  - zh-cn / de-de "Error Detail" cleared now that the handback is in sync (no more
#    stale-version warnings)
#  - Column widths widened/narrowed to fit the new text

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de status columns so the longer text fits.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Latest Handback DateTime (shared by both data rows) advances to the new run.
$zhcn.Range("K2").Value = "2016-08-31 07:34:59"
$zhcn.Range("K3").Value = "2016-08-31 07:34:59"

# Error Detail: the handback is now in sync, so the stale-version warnings clear.
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("K2").Value = "2016-08-31 07:35:26"
$dede.Range("K3").Value = "2016-08-31 07:35:26"

$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
